$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh the aggregate stats now that trade #37 has closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.59   # Current Capital
$summary.Range("B4").Value = -0.41     # Total P&L $
$summary.Range("B5").Value = -0.22     # Total P&L %
$summary.Range("B6").Value = 37        # Total Trades
$summary.Range("B7").Value = 12        # Winning Trades
$summary.Range("B9").Value = 32.43     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) reflects the same update.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.59      # Capital
$status.Range("D4").Value = 37         # Trades
$status.Range("E4").Value = -0.41      # P&L $
$status.Range("F4").Value = -0.41      # P&L %
$status.Range("G4").Value = 32.43      # Win Rate %

# ---------------------------------------------------------------------------
# Append trade #37 as a new row (row 38) to both the "All Trades" sheet and
# the per-strategy "MarketMaking" sheet.
# ---------------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Range("A38").Value = 37
    # B/C hold date- and time-looking text; force text formatting so Excel
    # doesn't coerce them into date/time serials, then restore the default
    # "Normal" style so the cell doesn't carry a stray number format.
    $ws.Range("B38:C38").NumberFormat = "@"
    $ws.Range("B38").Value = "2026-02-17"
    $ws.Range("C38").Value = "08:32:52"
    $ws.Range("B38:C38").Style = "Normal"
    $ws.Range("D38").Value = "MarketMaking"
    $ws.Range("E38").Value = "DOWN"
    $ws.Range("F38").Value = 0.64
    $ws.Range("G38").Value = 0.66
    $ws.Range("H38").Value = "CLOSED"
    $ws.Range("I38").Value = 3.125
    $ws.Range("J38").Value = 0.02
    $ws.Range("K38").Value = 99.59
    $ws.Range("L38").Value = 0
    $ws.Range("M38").Value = 0
    $ws.Range("N38").Value = 0.6
    $ws.Range("O38").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P38").Value = "early_exit"
    $ws.Range("Q38").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
